# Daily refresh of the cryptos price/volume table (GitHub Actions job).
# Price cells that look like a plain decimal number (e.g. "226.37") are
# written with a leading apostrophe so Excel keeps them as text (matching
# the source sheet, where every Price/Volume cell is a string, not a
# number) instead of silently coercing them to floating point.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.353.49'
$ws.Range("E2").Value = '  +0.77%  '
$ws.Range("D3").Value = '1.787.64'
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '''226.37'
$ws.Range("E5").Value = '  +0.41%  '
$ws.Range("D6").Value = '''0.557'
$ws.Range("E6").Value = '  +2.17%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("E8").Value = '  +3.54%  '
$ws.Range("E9").Value = '  +1.18%  '
$ws.Range("E10").Value = '  +0.36%  '
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("E12").Value = '  +0.39%  '
$ws.Range("D13").Value = '''11.22'
$ws.Range("E13").Value = '  +2.74%  '
$ws.Range("D14").Value = '1.795.66'
$ws.Range("E14").Value = '  +0.98%  '
$ws.Range("E15").Value = '  +2.28%  '
$ws.Range("D16").Value = '34.342.41'
$ws.Range("E16").Value = '  +0.76%  '
$ws.Range("E17").Value = '  +2.70%  '
$ws.Range("D18").Value = '''68.42'
$ws.Range("E18").Value = '  +1.22%  '
$ws.Range("D19").Value = '''244.77'
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("E20").Value = '  +0.75%  '
$ws.Range("D21").Value = '''11.25'
$ws.Range("E21").Value = '  +3.40%  '
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("E23").Value = '  +1.52%  '
$ws.Range("D24").Value = '''168.87'
$ws.Range("E24").Value = '  +4.46%  '
$ws.Range("D25").Value = '''2.07'
$ws.Range("E25").Value = '  +1.67%  '
$ws.Range("D26").Value = '''7.33'
$ws.Range("E26").Value = '  +3.32%  '
$ws.Range("D27").Value = '''16.52'
$ws.Range("E27").Value = '  +1.91%  '
$ws.Range("E28").Value = '  +1.64%  '
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("E30").Value = '  +9.43%  '
$ws.Range("E31").Value = '  +1.91%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '''3.79'
$ws.Range("E32").Value = '  +2.82%  '
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").Value = '''1.23'
$ws.Range("E33").Value = '  +0.11%  '
$ws.Range("E34").Value = '  +1.41%  '
$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D35").Value = '''2.58'
$ws.Range("E35").Value = '  +4.85%  '
$ws.Range("B36").Value = 'Maker'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D36").Value = '1.408.40'
$ws.Range("E36").Value = '  -2.54%  '
$ws.Range("D37").Value = '''0.682'
$ws.Range("E37").Value = '  +4.97%  '
$ws.Range("E38").Value = '  +3.09%  '
$ws.Range("E39").Value = '  +0.44%  '
$ws.Range("D40").Value = '''84.43'
$ws.Range("E40").Value = '  +5.14%  '
$ws.Range("D41").Value = '''2.40'
$ws.Range("E41").Value = '  +0.27%  '
$ws.Range("E42").Value = '  +2.43%  '
$ws.Range("D43").Value = '''0.938'
$ws.Range("E43").Value = '  +2.65%  '
$ws.Range("D44").Value = '''14.01'
$ws.Range("E44").Value = '  +2.59%  '
$ws.Range("E45").Value = '  +1.85%  '
$ws.Range("E46").Value = '  +2.56%  '
$ws.Range("D47").Value = '''6.07'
$ws.Range("E47").Value = '  +0.58%  '
$ws.Range("E48").Value = '  +0.39%  '
$ws.Range("D49").Value = '''105.29'
$ws.Range("E49").Value = '  +1.05%  '
$ws.Range("E50").Value = '  -0.11%  '
$ws.Range("E51").Value = '  -1.52%  '
